$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update x/y coordinate values for rows 2-50 (booth entries 1-49)
$ws.Cells.Item(2, 2).Value2 = 148
$ws.Cells.Item(2, 3).Value2 = 130.5
$ws.Cells.Item(3, 2).Value2 = 224
$ws.Cells.Item(3, 3).Value2 = 130.5
$ws.Cells.Item(4, 2).Value2 = 383
$ws.Cells.Item(4, 3).Value2 = 132
$ws.Cells.Item(5, 2).Value2 = 462
$ws.Cells.Item(5, 3).Value2 = 132
$ws.Cells.Item(6, 2).Value2 = 1183
$ws.Cells.Item(6, 3).Value2 = 136.5
$ws.Cells.Item(7, 2).Value2 = 629
$ws.Cells.Item(7, 3).Value2 = 137
$ws.Cells.Item(8, 2).Value2 = 703
$ws.Cells.Item(8, 3).Value2 = 137
$ws.Cells.Item(9, 2).Value2 = 1110
$ws.Cells.Item(9, 3).Value2 = 138
$ws.Cells.Item(10, 2).Value2 = 875
$ws.Cells.Item(10, 3).Value2 = 140
$ws.Cells.Item(11, 2).Value2 = 950
$ws.Cells.Item(11, 3).Value2 = 140
$ws.Cells.Item(12, 2).Value2 = 148
$ws.Cells.Item(12, 3).Value2 = 241.5
$ws.Cells.Item(13, 2).Value2 = 222
$ws.Cells.Item(13, 3).Value2 = 241.5
$ws.Cells.Item(14, 2).Value2 = 462
$ws.Cells.Item(14, 3).Value2 = 248
$ws.Cells.Item(15, 2).Value2 = 1110
$ws.Cells.Item(15, 3).Value2 = 249.5
$ws.Cells.Item(16, 2).Value2 = 1183
$ws.Cells.Item(16, 3).Value2 = 249.5
$ws.Cells.Item(17, 2).Value2 = 506
$ws.Cells.Item(17, 3).Value2 = 251.25
$ws.Cells.Item(18, 2).Value2 = 703
$ws.Cells.Item(18, 3).Value2 = 254.5
$ws.Cells.Item(19, 2).Value2 = 950
$ws.Cells.Item(19, 3).Value2 = 258
$ws.Cells.Item(20, 2).Value2 = 875
$ws.Cells.Item(20, 3).Value2 = 261.5
$ws.Cells.Item(21, 2).Value2 = 150
$ws.Cells.Item(21, 3).Value2 = 347.5
$ws.Cells.Item(22, 2).Value2 = 222
$ws.Cells.Item(22, 3).Value2 = 350
$ws.Cells.Item(23, 2).Value2 = 383
$ws.Cells.Item(23, 3).Value2 = 358
$ws.Cells.Item(24, 2).Value2 = 462
$ws.Cells.Item(24, 3).Value2 = 358
$ws.Cells.Item(25, 2).Value2 = 1110
$ws.Cells.Item(25, 3).Value2 = 358
$ws.Cells.Item(26, 2).Value2 = 1183
$ws.Cells.Item(26, 3).Value2 = 358
$ws.Cells.Item(27, 2).Value2 = 629
$ws.Cells.Item(27, 3).Value2 = 362.5
$ws.Cells.Item(28, 2).Value2 = 703
$ws.Cells.Item(28, 3).Value2 = 362.5
$ws.Cells.Item(29, 2).Value2 = 875
$ws.Cells.Item(29, 3).Value2 = 370
$ws.Cells.Item(30, 2).Value2 = 950
$ws.Cells.Item(30, 3).Value2 = 372
$ws.Cells.Item(31, 2).Value2 = 222
$ws.Cells.Item(31, 3).Value2 = 468.5
$ws.Cells.Item(32, 2).Value2 = 150
$ws.Cells.Item(32, 3).Value2 = 469.5
$ws.Cells.Item(33, 2).Value2 = 383
$ws.Cells.Item(33, 3).Value2 = 474
$ws.Cells.Item(34, 2).Value2 = 462
$ws.Cells.Item(34, 3).Value2 = 474
$ws.Cells.Item(35, 2).Value2 = 1110
$ws.Cells.Item(35, 3).Value2 = 474
$ws.Cells.Item(36, 2).Value2 = 1183
$ws.Cells.Item(36, 3).Value2 = 474
$ws.Cells.Item(37, 2).Value2 = 629
$ws.Cells.Item(37, 3).Value2 = 478.5
$ws.Cells.Item(38, 2).Value2 = 703
$ws.Cells.Item(38, 3).Value2 = 478.5
$ws.Cells.Item(39, 2).Value2 = 875
$ws.Cells.Item(39, 3).Value2 = 486
$ws.Cells.Item(40, 2).Value2 = 950
$ws.Cells.Item(40, 3).Value2 = 487
$ws.Cells.Item(41, 2).Value2 = 383
$ws.Cells.Item(41, 3).Value2 = 590
$ws.Cells.Item(42, 2).Value2 = 462
$ws.Cells.Item(42, 3).Value2 = 590
$ws.Cells.Item(43, 2).Value2 = 1110
$ws.Cells.Item(43, 3).Value2 = 590
$ws.Cells.Item(44, 2).Value2 = 1183
$ws.Cells.Item(44, 3).Value2 = 590
$ws.Cells.Item(45, 2).Value2 = 150
$ws.Cells.Item(45, 3).Value2 = 591.5
$ws.Cells.Item(46, 2).Value2 = 222
$ws.Cells.Item(46, 3).Value2 = 592.5
$ws.Cells.Item(47, 2).Value2 = 875
$ws.Cells.Item(47, 3).Value2 = 594
$ws.Cells.Item(48, 2).Value2 = 629
$ws.Cells.Item(48, 3).Value2 = 594.5
$ws.Cells.Item(49, 2).Value2 = 703
$ws.Cells.Item(49, 3).Value2 = 594.5
$ws.Cells.Item(50, 2).Value2 = 950
$ws.Cells.Item(50, 3).Value2 = 598

# Remove the now-obsolete trailing rows (previously rows 51-65, booths 50-64)
$ws.Range("A51:C65").EntireRow.Delete()
